$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Gdf9", "Bmpr2", "ECs", 3, 1, 0.9029553333333333, 2.708866, 0.2353208431348719, 0.2353208431348719, 3, 1, 40.70766766666667, 122.123003, 0.3776398983502007, 0.3776398983502007, 36.75720562717756, 330.814850644598, 0.08886653928113657, 0.08886653928113657),
    @("ECs", "Gdf9", "Bmpr2", "FAPs", 3, 1, 0.9029553333333333, 2.708866, 0.2353208431348719, 0.2353208431348719, 3, 1, 39.715023, 119.145069, 0.3684312589831062, 0.3684312589831062, 35.860891831306, 322.748026481754, 0.0866995545011469, 0.08669955450114691),
    @("ECs", "Gdf9", "Bmpr2", "sCs", 3, 1, 0.9029553333333333, 2.708866, 0.2353208431348719, 0.2353208431348719, 3, 1, 27.37224266666666, 82.116728, 0.253928842666693, 0.253928842666693, 24.71591250116089, 222.443212510448, 0.05975474935258844, 0.05975474935258845),
    @("FAPs", "Gdf9", "Bmpr2", "ECs", 3, 1, 1.446898, 4.340694, 0.3770787376970584, 0.3770787376970584, 3, 1, 40.70766766666667, 122.123003, 0.3776398983502007, 0.3776398983502007, 58.89984293156466, 530.098586384082, 0.1423999761739392, 0.1423999761739391),
    @("FAPs", "Gdf9", "Bmpr2", "FAPs", 3, 1, 1.446898, 4.340694, 0.3770787376970584, 0.3770787376970584, 3, 1, 39.715023, 119.145069, 0.3684312589831062, 0.3684312589831062, 57.463587348654, 517.1722861378861, 0.1389275940654877, 0.1389275940654877),
    @("FAPs", "Gdf9", "Bmpr2", "sCs", 3, 1, 1.446898, 4.340694, 0.3770787376970584, 0.3770787376970584, 3, 1, 27.37224266666666, 82.116728, 0.253928842666693, 0.253928842666693, 39.60484316991467, 356.443588529232, 0.09575116745763156, 0.09575116745763156),
    @("sCs", "Gdf9", "Bmpr2", "ECs", 2, 0.6666666666666666, 1.487271, 4.461812999999999, 0.3876004191680696, 0.3876004191680697, 3, 1, 40.70766766666667, 122.123003, 0.3776398983502007, 0.3776398983502007, 60.54333359827099, 544.890002384439, 0.146373382895125, 0.146373382895125),
    @("sCs", "Gdf9", "Bmpr2", "FAPs", 2, 0.6666666666666666, 1.487271, 4.461812999999999, 0.3876004191680696, 0.3876004191680697, 3, 1, 39.715023, 119.145069, 0.3684312589831062, 0.3684312589831062, 59.06700197223299, 531.6030177500969, 0.1428041104164716, 0.1428041104164716),
    @("sCs", "Gdf9", "Bmpr2", "sCs", 2, 0.6666666666666666, 1.487271, 4.461812999999999, 0.3876004191680696, 0.3876004191680697, 3, 1, 27.37224266666666, 82.116728, 0.253928842666693, 0.253928842666693, 40.70994272309599, 366.3894845078639, 0.09842292585647303, 0.09842292585647304),
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowVals[$c]
    }
}